$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (Oyuncu Adı / Pozisyon / Takım)
$data = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
